# Rod_test_analysis.xlsx edit: add confidence-interval columns and display/error-bar columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert two columns before J for char_strength 95% upper/lower bounds ---
$ws.Columns("J:K").Insert()
$ws.Range("J1").Value = "char_strength_95%_upper"
$ws.Range("K1").Value = "char_strength_95%_lower"

# --- Step 2: insert two columns before M (current position of "Design Strength")
#             for weibull_modulus 95% upper/lower bounds ---
$ws.Columns("M:N").Insert()
$ws.Range("M1").Value = "weibull_modulus_95%_upper"
$ws.Range("N1").Value = "weibull_modulus_95%_lower"

# --- Step 3: append two new columns at the end for display/error-bar strings,
#             copying the header style (bold, centered, bordered) from an
#             existing header cell ---
$ws.Range("A1").Copy()
$ws.Range("Y1:Z1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("Y1").Value = "char_strength_disp"
$ws.Range("Z1").Value = "weibull_modulus_disp"

# --- Step 4: fill in the confidence-interval values (rows with char_strength data) ---
$ws.Range("J2").Value = 301.1757812192096
$ws.Range("K2").Value = 196.2179041121504
$ws.Range("M2").Value = 6.33233774559996
$ws.Range("N2").Value = 1.893633397658162

$ws.Range("J3").Value = 145.8141782341235
$ws.Range("K3").Value = 123.4988439663269
$ws.Range("M3").Value = 12.99038824738473
$ws.Range("N3").Value = 4.716798374083251

$ws.Range("J4").Value = 176.0407780413132
$ws.Range("K4").Value = 123.8066111143724
$ws.Range("M4").Value = 6.257881243360847
$ws.Range("N4").Value = 2.226505991283302

$ws.Range("J5").Value = 91.56091174079913
$ws.Range("K5").Value = 71.81408476429806
$ws.Range("M5").Value = 19.78614227873543
$ws.Range("N5").Value = 3.806805126791093

$ws.Range("J8").Value = 83.39910513087202
$ws.Range("K8").Value = 50.78399434695405
$ws.Range("M8").Value = 10.70482092649276
$ws.Range("N8").Value = 1.834917442972895

# --- Step 5: fill in the display strings (mean with error-bar range beneath) ---
$ws.Range("Y2").Value = "243`n(301,196)"
$ws.Range("Z2").Value = "3.5`n(6.3,1.9)"

$ws.Range("Y3").Value = "134`n(146,123)"
$ws.Range("Z3").Value = "7.8`n(13.0,4.7)"

$ws.Range("Y4").Value = "148`n(176,124)"
$ws.Range("Z4").Value = "3.7`n(6.3,2.2)"

$ws.Range("Y5").Value = "81`n(92,72)"
$ws.Range("Z5").Value = "8.7`n(19.8,3.8)"

$ws.Range("Y8").Value = "65`n(83,51)"
$ws.Range("Z8").Value = "4.4`n(10.7,1.8)"

# Re-fit row heights so the multi-line display strings above don't leave
# behind an explicit/custom row height (keeps rows at their default height).
$ws.Rows("1:10").AutoFit()

Write-Host ("Done. UsedRange: " + $ws.UsedRange.Address())
